$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46049.01041666666, 0),
    @(3, 46049.02083333334, 0),
    @(4, 46049.03125, 0),
    @(5, 46049.04166666666, 0),
    @(6, 46049.05208333334, 0.29),
    @(7, 46049.0625, 0),
    @(8, 46049.07291666666, 0),
    @(9, 46049.08333333334, 0),
    @(10, 46049.09375, 0),
    @(11, 46049.10416666666, 0),
    @(12, 46049.11458333334, 0),
    @(13, 46049.125, 0),
    @(14, 46049.13541666666, 0.45),
    @(15, 46049.14583333334, 0),
    @(16, 46049.15625, 0),
    @(17, 46049.16666666666, 0),
    @(18, 46049.17708333334, 0),
    @(19, 46049.1875, 0),
    @(20, 46049.19791666666, 0),
    @(21, 46049.20833333334, 0),
    @(22, 46049.21875, 0.47),
    @(23, 46049.22916666666, 0.482),
    @(24, 46049.23958333334, 0.486),
    @(25, 46049.25, 0.496),
    @(26, 46049.26041666666, 1.191),
    @(27, 46049.27083333334, 1.586),
    @(28, 46049.28125, 2.382),
    @(29, 46049.29166666666, 4.439),
    @(30, 46049.30208333334, 10.765),
    @(31, 46049.3125, 17.464),
    @(32, 46049.32291666666, 27.565),
    @(33, 46049.33333333334, 39.337),
    @(34, 46049.34375, 67.52800000000001),
    @(35, 46049.35416666666, 82.163),
    @(36, 46049.36458333334, 98.011),
    @(37, 46049.375, 113.282),
    @(38, 46049.38541666666, 138.55),
    @(39, 46049.39583333334, 154.292),
    @(40, 46049.40625, 170.669),
    @(41, 46049.41666666666, 190.726),
    @(42, 46049.42708333334, 219.09),
    @(43, 46049.4375, 237.711),
    @(44, 46049.44791666666, 251.474),
    @(45, 46049.45833333334, 262.674),
    @(46, 46049.46875, 277.238),
    @(47, 46049.47916666666, 282.434),
    @(48, 46049.48958333334, 284.54),
    @(49, 46049.5, 284.12),
    @(50, 46049.51041666666, 281.541),
    @(51, 46049.52083333334, 277.104),
    @(52, 46049.53125, 269.716),
    @(53, 46049.54166666666, 255.187),
    @(54, 46049.55208333334, 229.187),
    @(55, 46049.5625, 212.311),
    @(56, 46049.57291666666, 200.873),
    @(57, 46049.58333333334, 186.273),
    @(58, 46049.59375, 157.361),
    @(59, 46049.60416666666, 140.129),
    @(60, 46049.61458333334, 121.163),
    @(61, 46049.625, 106.789),
    @(62, 46049.63541666666, 72.187),
    @(63, 46049.64583333334, 58.42),
    @(64, 46049.65625, 46.294),
    @(65, 46049.66666666666, 34.819),
    @(66, 46049.67708333334, 23.492),
    @(67, 46049.6875, 11.932),
    @(68, 46049.69791666666, 10.091),
    @(69, 46049.70833333334, 8.919),
    @(70, 46049.71875, 2.662),
    @(71, 46049.72916666666, 2.686),
    @(72, 46049.73958333334, 2.734),
    @(73, 46049.75, 2.802),
    @(74, 46049.76041666666, 2.65),
    @(75, 46049.77083333334, 0.65),
    @(76, 46049.78125, 0),
    @(77, 46049.79166666666, 0),
    @(78, 46049.80208333334, 0.49),
    @(79, 46049.8125, 0),
    @(80, 46049.82291666666, 2.49),
    @(81, 46049.83333333334, 0),
    @(82, 46049.84375, 2.65),
    @(83, 46049.85416666666, 0),
    @(84, 46049.86458333334, 0),
    @(85, 46049.875, 0),
    @(86, 46049.88541666666, 2.45),
    @(87, 46049.89583333334, 0),
    @(88, 46049.90625, 0.45),
    @(89, 46049.91666666666, 0),
    @(90, 46049.92708333334, 0),
    @(91, 46049.9375, 0),
    @(92, 46049.94791666666, 0),
    @(93, 46049.95833333334, 0),
    @(94, 46049.96875, 0),
    @(95, 46049.97916666666, 0),
    @(96, 46049.98958333334, 0),
    @(97, 46050.0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $aVal = $row[1]
    $bVal = $row[2]
    $ws.Cells.Item($r, 1).Value = $aVal
    $ws.Cells.Item($r, 2).Value = $bVal
}
